$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear any special formatting on A6 so it reverts to the default style
$ws.Range("A6").ClearFormats()
$ws.Range("A6").Value = "locked_out_user"
$ws.Range("B6").Value = "secret_sauce"

# Add new rows of data
$ws.Range("B7").Value = "secret_sauce"
$ws.Range("A8").Value = "standard_user"
$ws.Range("A9").Value = "abc"
$ws.Range("B9").Value = 123

$ws.Range("B8").Select()
